$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content updates (done before the column delete, while columns are still
#     C..I, so the affected addresses match the *current* layout) ---

# Competitor list: "Quandoo" -> "BookingNinja"
$ws.Range("G8").Value2 = "BookingNinja"

# "Most important findings" box for Strengths (I9, later becomes H9):
# replace the old sleek-design note with the new social-media note.
$ws.Range("I9").Value2 = "Social media connectivity" + [char]10 + "Intuitive for the customer - make it party size, date, then times available." + [char]10 + "Appointment reminders sent to the user via email or SMS" + [char]10

# New notes added for "BookingNinja" column (F) across the Strengths rows
# and one Weaknesses row.
$ws.Range("F9").Value2 = "Provides the business with its own branded application that the user can download and use to book their services."
$ws.Range("F10").Value2 = "Connectivity with social platforms making booking a lot easier."
$ws.Range("F11").Value2 = "The website is completely personalised for the business rather than being a plugin. This means that the website can be a seamless booking experience for the user."
$ws.Range("F12").Value2 = "Reminders are scheduled for the user to alert them when their booking is. Can be sent via email or SMS. This will decrease the chance of people missing their booking slot."
$ws.Range("F14").Value2 = "Doesn't allow the user to choose how many people the booking is for "

# --- Structural change: remove the empty filler column H so the
#     "Most important findings" column slides left from I to H ---
$ws.Columns.Item(8).Delete()

# --- Row heights that grew to fit the newly-added text ---
$ws.Rows.Item(10).RowHeight = 107.25
$ws.Rows.Item(14).RowHeight = 122.25

# --- View state tweaks ---
$ws.Application.ActiveWindow.Zoom = 100
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("G10").Select()
